# "lock in current version" -------------------------------------------
# Rename the four response-image headers, add a new trailing "iti"
# (inter-trial-interval) column with its per-trial latency values,
# refresh the cue word list (column C) with a new set of German verbs,
# and repoint the "correct" column (F) at the buffered copies of the
# target images.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row: rename image1..image4 -> image_top/right/bottom/left,
#      and introduce the new "iti" header in column L --------------------
$ws.Range("H1").Value = "image_top"
$ws.Range("I1").Value = "image_right"
$ws.Range("J1").Value = "image_bottom"
$ws.Range("K1").Value = "image_left"
$ws.Range("L1").Value = "iti"

# New header cell picks up the same bold / centered-top look as the rest
# of row 1, plus a thin left/right border to set it apart as a new block.
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("L1").VerticalAlignment = -4160     # xlTop
$ws.Range("L1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("L1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

# ---- Column C ("word"): new cue-word set --------------------------------
$ws.Range("C2").Value = "anrufen"
$ws.Range("C3").Value = "verkaufen"
$ws.Range("C4").Value = "verletzen"
$ws.Range("C5").Value = "heiraten"
$ws.Range("C6").Value = "gefallen"
$ws.Range("C7").Value = "gewinnen"
$ws.Range("C8").Value = "bestimmen"
$ws.Range("C9").Value = "versprechen"

# ---- Column F ("correct"): point at the buffered image copies ----------
$ws.Range("F2").Value = "buffer/dog/dog275.png"
$ws.Range("F3").Value = "buffer/house/house270.png"
$ws.Range("F4").Value = "buffer/flower/flower270.png"
$ws.Range("F5").Value = "buffer/face/face215.png"
$ws.Range("F6").Value = "buffer/dog/dog276.png"
$ws.Range("F7").Value = "buffer/house/house271.png"
$ws.Range("F8").Value = "buffer/face/face271.png"
$ws.Range("F9").Value = "buffer/flower/flower271.png"

# ---- Column L ("iti"): new per-trial inter-trial-interval values --------
$ws.Range("L2").Value = 1.015820934783024
$ws.Range("L3").Value = 1.0716902907541741
$ws.Range("L4").Value = 1.047075324721138
$ws.Range("L5").Value = 1.1351010053757959
$ws.Range("L6").Value = 1.101215866370167
$ws.Range("L7").Value = 1.232576476035504
$ws.Range("L8").Value = 1.236406268413865
$ws.Range("L9").Value = 1.0490166754144501

# ---- Match the saved cursor position recorded in the workbook ----------
$ws.Range("F9").Select() | Out-Null
